$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.518.93"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "2.643.97"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'603.07"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "'156.38"
$ws.Range("E6").Value = "  +2.58%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "2.642.34"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").Value = "'0.123"
$ws.Range("E10").Value = "  +7.10%  "
$ws.Range("D11").Value = "'0.401"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "'29.66"
$ws.Range("E14").Value = "  +5.39%  "
$ws.Range("D15").Value = "'0.0000193"
$ws.Range("E15").Value = "  +12.56%  "
$ws.Range("D16").Value = "3.119.45"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "65.231.84"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "2.655.54"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").Value = "'4.85"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "'356.83"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").Value = "'7.40"
$ws.Range("E22").Value = "  +3.78%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'69.51"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("E26").Value = "  +14.62%  "
$ws.Range("D27").Value = "'9.36"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").Value = "'1.62"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").Value = "'8.09"
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  +4.27%  "
$ws.Range("D33").Value = "'529.82"
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("D34").Value = "'1.77"
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("D35").Value = "'5.52"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").Value = "'0.431"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").Value = "'20.57"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("D39").Value = "'161.59"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.97"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("D44").Value = "'163.97"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "'4.13"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'2.35"
$ws.Range("E46").Value = "  +5.93%  "
$ws.Range("D47").Value = "'0.0606"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("D48").Value = "'22.77"
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("D49").Value = "'0.650"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").Value = "'0.0980"
$ws.Range("E51").Value = "  +0.61%  "
